$wb = $excel.ActiveWorkbook

# --- Sheet "Prix Spot": add column O (28-jun) ---
$ws1 = $wb.Worksheets.Item("Prix Spot")

$ws1.Range("N1").Copy()
$ws1.Range("O1").PasteSpecial(-4122) # xlPasteFormats - match header style of the other date columns
$ws1.Range("O1").Value = "28-jun"

$o1Values = @(95.11, 78.5, 65.56999999999999, 35.07, 51.29, 44.3, 31.81, 47.38, 27.8, 2.34, 0.01, 0, -0.01, -0.02, -0.02, -0.02, 2.12, 10.26, 65, 96.26000000000001, 110, 105.09, 117.48, 103.5)

for ($i = 0; $i -lt $o1Values.Length; $i++) {
    $row = $i + 2
    $ws1.Cells.Item($row, 15).Value = $o1Values[$i]
}

# --- Sheet "Gaz": add row 12 (2025-06-26, 32.625) ---
$ws2 = $wb.Worksheets.Item("Gaz")
# Force text so "2025-06-26" isn't auto-converted to a date serial, then
# reset the cell style back to Normal so no stray number format sticks.
$ws2.Range("A12").NumberFormat = "@"
$ws2.Range("A12").Value = "2025-06-26"
$ws2.Range("A12").Style = "Normal"
$ws2.Range("B12").Value = 32.625

# --- Sheet "CO2": add row 12 (2025-06-26, 69.45999999999999) ---
$ws3 = $wb.Worksheets.Item("CO2")
$ws3.Range("A12").NumberFormat = "@"
$ws3.Range("A12").Value = "2025-06-26"
$ws3.Range("A12").Style = "Normal"
$ws3.Range("B12").Value = 69.45999999999999
